$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values: B1,C1,D1,E1 updated to 15,16,15,16
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON): B2 and D2 updated; C2 and E2 cleared entirely
$ws.Range("B2").Value = 10.058760961894023
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 11.217134232315388
$ws.Range("E2").ClearContents()

# Row 3 (STR): B3,C3,D3,E3 updated
$ws.Range("B3").Value = 8.6204950727559577
$ws.Range("C3").Value = -3.6976792365190221
$ws.Range("D3").Value = 8.4310905856397014
$ws.Range("E3").Value = -6.8681887904253571

# Update the selection to match the new used range B1:E3
$ws.Range("B1:E3").Select()
